$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename foils from batch "2" to batch "3"
$ws.Range("A2").Value = "Zr3"
$ws.Range("A3").Value = "In3"
$ws.Range("A4").Value = "Ni3"
$ws.Range("A5").Value = "Au3"
$ws.Range("A8").Value = "Al3"

# Updated measurements
$ws.Range("E3").Value = 0.01
$ws.Range("D5").Value = 0.1

# New uncertainty column (I) with propagated-error formula
$ws.Range("I2").Formula = "=SQRT((C2/B2)^2+(E2/D2)^2+(G2/F2)^2)*H2"
$ws.Range("I3:I8").Formula = "=SQRT((C3/B3)^2+(E3/D3)^2+(G3/F3)^2)*H3"

# Update the saved selection to match the source workbook
$ws.Range("E5").Select()

$wb.Save()
